function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "BaseConfig": insert a new "partition_num" row right after
# "dis_channel" (before "dws_url"), shifting the rest of the rows down.
# ---------------------------------------------------------------------------
$base = $wb.Worksheets.Item("BaseConfig")
$base.Rows.Item(3).Insert()
$base.Range("A3").Value = "partition_num"
$base.Range("B3").Value = "2"

# ---------------------------------------------------------------------------
# Sheet "TableConfig": the per-table "partition" column is gone; instead the
# "channel" column now lives in column B (previously a trailing, largely
# empty column) and two brand-new columns are added describing, per
# partition, whether the target needs to be created and which app group it
# belongs to.
# ---------------------------------------------------------------------------
$tbl = $wb.Worksheets.Item("TableConfig")

# Drop the old numeric "partition" data (column B); col_map/source_column
# stay untouched in columns C/D.
$tbl.Range("B2:B5").Clear()

# Column B header becomes "channel" (used to be the header of column E).
$tbl.Range("B1").Value = "channel"

# New columns: need_create (E) and group_name (F).
$tbl.Range("E1").Value = "need_create"
$tbl.Range("F1").Value = "group_name"

$tbl.Range("E2").Value = "y"
$tbl.Range("F2").Value = "app1"

$tbl.Range("E3").Value = "y"
$tbl.Range("F3").Value = "app2"

$tbl.Range("E4").Value = "n"
$tbl.Range("F4").Value = "app1"

$tbl.Range("E5").Value = "y"
$tbl.Range("F5").Value = "app2"

# Match the original column styling used by the header row ("channel" header
# keeps the same highlighted look, just without the loud yellow fill).
$tbl.Range("B1").Style = $tbl.Range("A1").Style
$tbl.Range("E1:F1").Style = $tbl.Range("A1").Style
$tbl.Range("A2:A5,C2:C5,D2:D5").Style = $tbl.Range("A2").Style

# The "channel" header cell keeps its own highlight fill, just recolored
# from yellow to a plain white background.
$tbl.Range("B1").Interior.Color = RGB(255, 255, 255)

# Column widths: drop the old dedicated "partition" column width, widen the
# source_column column, and size the two new columns.
$tbl.Columns.Item(4).ColumnWidth = 50.26953125
$tbl.Columns.Item(5).ColumnWidth = 12.453125
$tbl.Columns.Item(6).ColumnWidth = 11.36328125

# Keep the data-entry rows (3-5) for the new columns free of the header
# style like the rest of the sheet's body rows below row 2.
$tbl.Range("E3,F3,E4,F4,E5,F5").Style = "Normal"
